# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.6545652718822623, 0.3048912486333797, 3.223369029078222, 0.5333859586016987, 4.716211508195562)
    3 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    4 = @(0.6545652718822623, 0.04103571897497393, 0.7210945179870265, 0.5333859586016987, 1.950081467445961)
    5 = @(3.272327238179451, 1.626987699542094, 18.71679738969934, 13.86384647080068, 37.47995879822157)
    6 = @(0.6545652718822623, 1.626987699542094, 0.1496068669990043, 13.86384647080068, 16.29500630922404)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E - IP
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G - sum
}

$wb.Save()
